$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Species GBIF", "Species Original"),
    @("Acer campestre", "Acer campestre"),
    @("Acer negundo", "Acer negundo"),
    @("Acer platanoides", "Acer platanoides"),
    @("Acer pseudoplatanus", "Acer pseudoplatanus"),
    @("Acer", "Acer sp."),
    @("Achillea millefolium", "Achillea millefolium"),
    @("Agrimonia eupatoria", "Agrimonia eupatoria"),
    @("Agrostis capillaris", "Agrostis capillaris"),
    @("Allium", "Allium sp."),
    @("Allium vineale", "Allium vineale"),
    @("Alopecurus pratensis", "Alopecurus pratensis"),
    @("Lysimachia arvensis", "Anagallis arvensis"),
    @("Anthemis arvensis", "Anthemis arvensis"),
    @("Anthoxanthum odoratum", "Anthoxanthum odoratum"),
    @("Anthriscus sylvestris", "Anthriscus sylvestris"),
    @("Apiaceae", "Apiaceae sp."),
    @("Arabidopsis thaliana", "Arabidopsis thaliana"),
    @("Arenaria serpyllifolia", "Arenaria serpyllifolia"),
    @("Arrhenatherum elatius", "Arrhenatherum elatius"),
    @("Artemisia vulgaris", "Artemisia vulgaris"),
    @("Asteraceae", "Asteraceae sp."),
    @("Avena sativa", "Avena sativa"),
    @("Bellis perennis", "Bellis perennis"),
    @("Betonica officinalis", "Betonica officinalis"),
    @("Betula", "Betula sp."),
    @("Borago officinalis", "Borago officinalis"),
    @("Brachypodium pinnatum", "Brachypodium pinnatum"),
    @("Brassica napus", "Brassica napus"),
    @("Brassica oleracea", "Brassica oleracea"),
    @("Brassicaceae", "Brassicaceae sp."),
    @("Bromus erectus", "Bromus erectus"),
    @("Bromus hordeaceus", "Bromus hordeaceus"),
    @("Bromus sterilis", "Bromus sterilis"),
    @("Calamagrostis epigejos", "Calamagrostis epigejos"),
    @("Campanula patula", "Campanula patula"),
    @("Campanula rapunculoides", "Campanula rapunculoides"),
    @("Campanula rotundifolia", "Campanula rotundifolia"),
    @("Capsella bursa-pastoris", "Capsella bursa-pastoris"),
    @("Carduus acanthoides", "Carduus acanthoides"),
    @("Carduus crispus", "Carduus crispus"),
    @("Carum carvi", "Carum carvi"),
    @("Centaurea jacea", "Centaurea jacea"),
    @("Centaurea scabiosa", "Centaurea scabiosa"),
    @("Cerastium glomeratum", "Cerastium glomeratum"),
    @("Cerastium holosteoides", "Cerastium holosteoides"),
    @("Chenopodium album", "Chenopodium album"),
    @("Asteraceae", "Chrysanthemum vulgare"),
    @("Cichorium intybus", "Cichorium intybus"),
    @("Cirsium arvense", "Cirsium arvense"),
    @("Cirsium vulgare", "Cirsium vulgare"),
    @("Clematis vitalba", "Clematis vitalba"),
    @("Colchicum autumnale", "Colchicum autumnale"),
    @("Convolvulus arvensis", "Convolvulus arvensis"),
    @("Erigeron canadensis", "Conyza canadensis"),
    @("Cornus sanguinea", "Cornus sanguinea"),
    @("Crataegus monogyna", "Crataegus monogyna"),
    @("Crepis biennis", "Crepis biennis"),
    @("Crepis capillaris", "Crepis capillaris"),
    @("Cynosurus cristatus", "Cynosurus cristatus"),
    @("Dactylis glomerata", "Dactylis glomerata"),
    @("Daucus carota", "Daucus carota"),
    @("Descurainia sophia", "Descurainia sophia"),
    @("Dianthus carthusianorum", "Dianthus carthusianorum"),
    @("Dianthus seguieri", "Dianthus seguieri"),
    @("Elymus repens", "Elymus repens"),
    @("Epilobium tetragonum", "Epilobium tetragonum"),
    @("Equisetum arvense", "Equisetum arvense"),
    @("Erigeron acris", "Erigeron acris"),
    @("Erigeron annuus", "Erigeron annuus"),
    @("Erodium cicutarium", "Erodium cicutarium"),
    @("Draba", "Erophila verna"),
    @("Eryngium campestre", "Eryngium campestre"),
    @("Falcaria vulgaris", "Falcaria vulgaris"),
    @("Fallopia convolvulus", "Fallopia convolvulus"),
    @("Poaceae", "Festuca arundinacea"),
    @("Festuca ovina", "Festuca ovina"),
    @("Poaceae", "Festuca pratensis"),
    @("Festuca rubra", "Festuca rubra"),
    @("Festuca rupicola", "Festuca rupicola"),
    @("Festuca", "Festuca sp."),
    @("Festulolium", "Festulolium"),
    @("Filipendula ulmaria", "Filipendula ulmaria"),
    @("Filipendula vulgaris", "Filipendula vulgaris"),
    @("Fraxinus excelsior", "Fraxinus excelsior"),
    @("Fumaria officinalis", "Fumaria officinalis"),
    @("Galium album", "Galium album"),
    @("Galium aparine", "Galium aparine"),
    @("Galium mollugo", "Galium mollugo"),
    @("Galium pomeranicum", "Galium x pomeranicum"),
    @("Galium verum", "Galium verum"),
    @("Geranium molle", "Geranium molle"),
    @("Geranium pratense", "Geranium pratense"),
    @("Geranium pusillum", "Geranium pusillum"),
    @("Geranium pyrenaicum", "Geranium pyrenaicum"),
    @("Geranium rotundifolium", "Geranium rotundifolium"),
    @("Geranium", "Geranium sp."),
    @("Geum urbanum", "Geum urbanum"),
    @("Glechoma hederacea", "Glechoma hederacea"),
    @("Helictotrichon pratense", "Helictotrichon pratense"),
    @("Avenula pubescens", "Helictotrichon pubescens"),
    @("Heracleum sphondylium", "Heracleum sphondylium"),
    @("Pilosella officinarum", "Hieracium pilosella"),
    @("Holcus lanatus", "Holcus lanatus"),
    @("Hypericum perforatum", "Hypericum perforatum"),
    @("Hypochaeris radicata", "Hypochaeris radicata"),
    @("Knautia arvensis", "Knautia arvensis"),
    @("Koeleria macrantha", "Koeleria macrantha"),
    @("Lactuca serriola", "Lactuca serriola"),
    @("Lamium amplexicaule", "Lamium amplexicaule"),
    @("Lamium purpureum", "Lamium purpureum"),
    @("Lapsana communis", "Lapsana communis"),
    @("Lathyrus pratensis", "Lathyrus pratensis"),
    @("Lathyrus tuberosus", "Lathyrus tuberosus"),
    @("Leontodon autumnalis", "Leontodon autumnalis"),
    @("Leontodon hispidus", "Leontodon hispidus"),
    @("Leucanthemum vulgare", "Leucanthemum vulgare"),
    @("Seseli libanotis", "Libanotis pyrenaica"),
    @("Linaria vulgaris", "Linaria vulgaris"),
    @("", ""),
    @("Lolium multiflorum", "Lolium multiflorum"),
    @("Lolium perenne", "Lolium perenne"),
    @("Lotus corniculatus", "Lotus corniculatus"),
    @("Lotus pedunculatus", "Lotus pedunculatus"),
    @("Silene flos-cuculi", "Lychnis flos-cuculi"),
    @("Lysimachia nummularia", "Lysimachia nummularia"),
    @("Malva sylvestris", "Malva sylvestris"),
    @("Asteraceae", "Matricaria inodora"),
    @("Medicago falcata", "Medicago falcata"),
    @("Medicago lupulina", "Medicago lupulina"),
    @("Medicago", "Medicago sp."),
    @("Medicago varia", "Medicago x varia"),
    @("Melampyrum pratense", "Melampyrum pratense"),
    @("Noccaea perfoliata", "Microthlaspi perfoliatum"),
    @("", ""),
    @("Myosotis arvensis", "Myosotis arvensis"),
    @("Papaver rhoeas", "Papaver rhoeas"),
    @("Pastinaca sativa", "Pastinaca sativa"),
    @("Petrorhagia prolifera", "Petrorhagia prolifera"),
    @("Phleum pratense", "Phleum pratense"),
    @("Picris hieracioides", "Picris hieracioides"),
    @("Pimpinella saxifraga", "Pimpinella saxifraga"),
    @("Plantago lanceolata", "Plantago lanceolata"),
    @("Plantago major", "Plantago major"),
    @("Plantago media", "Plantago media"),
    @("Poa angustifolia", "Poa angustifolia"),
    @("Poa annua", "Poa annua"),
    @("Poa bulbosa", "Poa bulbosa"),
    @("Poa pratensis", "Poa pratensis"),
    @("Poa trivialis", "Poa trivialis"),
    @("Poaceae", "Poaceae sp."),
    @("Polygonum aviculare", "Polygonum aviculare"),
    @("Populus canadensis", "Populus x canadensis"),
    @("Potentilla argentea", "Potentilla argentea"),
    @("Prunella vulgaris", "Prunella vulgaris"),
    @("Prunus avium", "Prunus avium"),
    @("Prunus mahaleb", "Prunus mahaleb"),
    @("Prunus", "Prunus sp."),
    @("Ranunculus acris", "Ranunculus acris"),
    @("Ranunculus auricomus", "Ranunculus auricomus"),
    @("Ranunculus bulbosus", "Ranunculus bulbosus"),
    @("Ranunculus repens", "Ranunculus repens"),
    @("", ""),
    @("Rhinanthus minor", "Rhinanthus minor"),
    @("Rubus caesius", "Rubus caesius"),
    @("Rubus idaeus", "Rubus idaeus"),
    @("Rubus", "Rubus sp."),
    @("Rumex crispus", "Rumex crispus"),
    @("Salvia pratensis", "Salvia pratensis"),
    @("Sambucus nigra", "Sambucus nigra"),
    @("Poterium sanguisorba", "Sanguisorba minor"),
    @("Saxifraga granulata", "Saxifraga granulata"),
    @("Scabiosa ochroleuca", "Scabiosa ochroleuca"),
    @("Coronilla varia", "Securigera varia"),
    @("", ""),
    @("Senecio jacobaea", "Senecio jacobaea"),
    @("Senecio", "Senecio sp."),
    @("Senecio vernalis", "Senecio vernalis"),
    @("Senecio vulgaris", "Senecio vulgaris"),
    @("Setaria viridis", "Setaria viridis"),
    @("Silaum silaus", "Silaum silaus"),
    @("Silene latifolia", "Silene latifolia"),
    @("Silene noctiflora", "Silene noctiflora"),
    @("Silene nutans", "Silene nutans"),
    @("Silene vulgaris", "Silene vulgaris"),
    @("Sinapis arvensis", "Sinapis arvensis"),
    @("Solanum tuberosum", "Solanum tuberosum"),
    @("Solidago canadensis", "Solidago canadensis"),
    @("Sonchus arvensis", "Sonchus arvensis"),
    @("Sonchus asper", "Sonchus asper"),
    @("Sonchus oleraceus", "Sonchus oleraceus"),
    @("Stachys recta", "Stachys recta"),
    @("Stellaria graminea", "Stellaria graminea"),
    @("Stellaria media", "Stellaria media"),
    @("Tanacetum vulgare", "Tanacetum vulgare"),
    @("Taraxacum officinale", "Taraxacum officinale"),
    @("Thlaspi arvense", "Thlaspi arvense"),
    @("Tragopogon dubius", "Tragopogon dubius"),
    @("Tragopogon orientalis", "Tragopogon orientalis"),
    @("Tragopogon pratensis", "Tragopogon pratensis"),
    @("Trifolium arvense", "Trifolium arvense"),
    @("Trifolium campestre", "Trifolium campestre"),
    @("Trifolium dubium", "Trifolium dubium"),
    @("Trifolium hybridum", "Trifolium hybridum"),
    @("Trifolium pratense", "Trifolium pratense"),
    @("Trifolium repens", "Trifolium repens"),
    @("Tripleurospermum inodorum", "Tripleurospermum perforatum"),
    @("Trisetum flavescens", "Trisetum flavescens"),
    @("Triticum", "Triticum sp."),
    @("Tussilago farfara", "Tussilago farfara"),
    @("Urtica urens", "Urtica urens"),
    @("Verbascum lychnitis", "Verbascum lychnitis"),
    @("Verbascum thapsus", "Verbascum thapsus"),
    @("Veronica agrestis", "Veronica agrestis"),
    @("Veronica arvensis", "Veronica arvensis"),
    @("Veronica chamaedrys", "Veronica chamaedrys"),
    @("Veronica hederifolia", "Veronica hederifolia"),
    @("Veronica persica", "Veronica persica"),
    @("", ""),
    @("Veronica serpyllifolia", "Veronica serpyllifolia"),
    @("Vicia", "Vicia angustifolia"),
    @("Vicia cracca", "Vicia cracca"),
    @("Vicia hirsuta", "Vicia hirsuta"),
    @("Vicia sativa", "Vicia sativa"),
    @("Vicia sepium", "Vicia sepium"),
    @("Vicia tetrasperma", "Vicia tetrasperma"),
    @("Vicia villosa", "Vicia villosa"),
    @("Viola arvensis", "Viola arvensis"),
    @("Festuca", "Vulpia myuros")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}